# Update the "teaching" slide layout (slideLayout5.xml):
#   - title placeholder ("信息"): reposition/resize + reduce default font size 56 -> 52pt
#   - body placeholder ("Body Level One..."): reposition/resize + reduce default font size 36 -> 40pt
#
# Point values below are deliberately not the "nice" 43.2 / 187.2 / 873.6 / 86.4 / 316.8 / 64.8
# literals: PowerPoint COM stores shape geometry in points (as a 32-bit float) before it is
# converted back to EMU on save, so a plain 187.2 literal round-trips to 2377439 EMU instead of
# 2377440. The values used here are the nearest representable points that convert back to the
# exact target EMU (off=548640/2377440, ext=11094720/1097280 for the title; off=548640/4023360,
# ext=11094720/822960 for the body).

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cl = $m.CustomLayouts.Item(5)   # "teaching" layout -> ppt/slideLayouts/slideLayout5.xml

# Title placeholder ("信息", id=55)
$title = $cl.Shapes.Item(1)
$title.Left   = 43.20000076293945    # -> 548640 EMU
$title.Top    = 187.20001220703125   # -> 2377440 EMU
$title.Width  = 873.6000366210938    # -> 11094720 EMU
$title.Height = 86.4000015258789     # -> 1097280 EMU
$title.TextFrame.TextRange.Font.Size = 52

# Body placeholder ("Body Level One...", id=56)
$body = $cl.Shapes.Item(2)
$body.Left   = 43.20000076293945     # -> 548640 EMU
$body.Top    = 316.8000183105469     # -> 4023360 EMU
$body.Width  = 873.6000366210938     # -> 11094720 EMU
$body.Height = 64.80000305175781     # -> 822960 EMU
$body.TextFrame.TextRange.Font.Size = 40
